# Regenerate the "K" column (G) values for the save_data sheet.
# These replace the old "Strike#" derived values (column G, rows 2-65)
# with the newly computed K values from the regenerated std/mean & s_vals
# calculation pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,1,0,1,1,1,0,0,1,0,1,1,0,1,2,0,0,4,2,1,4,2,2,0,1,1,0,0,0,2,1,1,1,1,0,2,1,1,3,0,1,3,1,0,1,2,1,2,2,1,1,1,5,1,1,1,2,1,1,2,2,1,2,0)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
